# C5-PowerPoint.pptx edit
# 1) Slide 6's table switches to a different table style.
# 2) The deck's theme colour palette is swapped from the "Integral" palette
#    to the standard Office palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Update the table style on the table in slide 6 ---------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{58F98BA0-76D4-4F7E-ABCB-E7BD8F6D9C2B}")
    }
}

# --- 2. Swap the theme colour scheme (Integral -> Office) ------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
$themeColors.Item(1).RGB  = 0          # dk1      000000
$themeColors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388    # dk2      44546A
$themeColors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407      # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Item(10).RGB = 4697456    # accent6  70AD47
$themeColors.Item(11).RGB = 12673797   # hlink    0563C1
$themeColors.Item(12).RGB = 7491477    # folHlink 954F72
